# Updates cryptos list data (price + 1h volume) cell values on Sheet1.
# Numeric-looking "Price" strings are written with a leading apostrophe
# (forces Excel to keep them as literal text, matching the source data
# which stores every Price/Volume value as text, e.g. "3.563.25").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.225.38'
$ws.Range('E2').Value = '  +1.27%  '
$ws.Range('D3').Value = '3.563.25'
$ws.Range('E3').Value = '  +1.95%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = "'619.60"
$ws.Range('E5').Value = '  +3.10%  '
$ws.Range('D6').Value = "'154.89"
$ws.Range('E6').Value = '  +3.83%  '
$ws.Range('D7').Value = '3.562.30'
$ws.Range('E7').Value = '  +1.93%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = "'0.492"
$ws.Range('E9').Value = '  +2.39%  '
$ws.Range('E10').Value = '  +5.73%  '
$ws.Range('D11').Value = "'7.46"
$ws.Range('E11').Value = '  +7.60%  '
$ws.Range('D12').Value = "'0.438"
$ws.Range('E12').Value = '  +3.95%  '
$ws.Range('B13').Value = 'ShibaInu'
$ws.Range('C13').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D13').Value = "'0.0000222"
$ws.Range('E13').Value = '  +2.64%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').Value = "'33.23"
$ws.Range('E14').Value = '  +5.60%  '
$ws.Range('D15').Value = '4.166.47'
$ws.Range('E15').Value = '  +1.97%  '
$ws.Range('D16').Value = '3.558.33'
$ws.Range('E16').Value = '  +1.83%  '
$ws.Range('D17').Value = '68.275.41'
$ws.Range('E17').Value = '  +1.48%  '
$ws.Range('E18').Value = '  -0.03%  '
$ws.Range('E19').Value = '  +5.55%  '
$ws.Range('D20').Value = "'15.99"
$ws.Range('E20').Value = '  +6.85%  '
$ws.Range('D21').Value = "'10.05"
$ws.Range('E21').Value = '  +11.81%  '
$ws.Range('D22').Value = "'454.60"
$ws.Range('E22').Value = '  +1.87%  '
$ws.Range('E23').Value = '  +4.44%  '
$ws.Range('D24').Value = "'78.48"
$ws.Range('E24').Value = '  +1.49%  '
$ws.Range('D25').Value = "'0.0000132"
$ws.Range('E25').Value = '  +2.96%  '
$ws.Range('D26').Value = '3.706.65'
$ws.Range('E26').Value = '  +1.98%  '
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('D28').Value = "'9.26"
$ws.Range('E28').Value = '  +13.43%  '
$ws.Range('D29').Value = "'10.52"
$ws.Range('E29').Value = '  +4.25%  '
$ws.Range('E30').Value = '  +12.02%  '
$ws.Range('E31').Value = '  +3.70%  '
$ws.Range('D32').Value = "'0.170"
$ws.Range('E32').Value = '  +3.83%  '
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('D34').Value = "'6.39"
$ws.Range('E34').Value = '  +5.64%  '
$ws.Range('D35').Value = "'26.14"
$ws.Range('E35').Value = '  +1.93%  '
$ws.Range('E36').Value = '  +4.95%  '
$ws.Range('D37').Value = '3.557.66'
$ws.Range('E37').Value = '  +2.07%  '
$ws.Range('D38').Value = "'8.24"
$ws.Range('E38').Value = '  +3.58%  '
$ws.Range('D39').Value = "'2.39"
$ws.Range('E39').Value = '  +8.85%  '
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('D41').Value = "'181.37"
$ws.Range('E41').Value = '  +3.77%  '
$ws.Range('D42').Value = "'0.0919"
$ws.Range('E42').Value = '  +5.35%  '
$ws.Range('E43').Value = '  -0.04%  '
$ws.Range('D44').Value = "'5.64"
$ws.Range('E44').Value = '  +4.86%  '
$ws.Range('D45').Value = "'31.17"
$ws.Range('E45').Value = '  +14.35%  '
$ws.Range('E46').Value = '  +2.24%  '
$ws.Range('D47').Value = "'46.18"
$ws.Range('E47').Value = '  +1.65%  '
$ws.Range('D48').Value = "'1.34"
$ws.Range('E48').Value = '  +5.59%  '
$ws.Range('D49').Value = "'2.65"
$ws.Range('E49').Value = '  +4.53%  '
$ws.Range('E50').Value = '  +3.66%  '
$ws.Range('E51').Value = '  +8.04%  '
